$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.049528901042879
$rowBF[0,2] = 1.054136975441201
$rowBF[0,3] = 1.053202837223017
$rowBF[0,4] = 1.063417280071903
$ws.Range("B2:F2").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.046267377102742
$rowIN[0,1] = 1.054566867383135
$rowIN[0,2] = 1.056881015881968
$rowIN[0,3] = 1.055949455611956
$rowIN[0,4] = 1.066135975885871
$rowIN[0,5] = 1.021951551960112
$ws.Range("I2:N2").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.050589452302866
$rowBF[0,2] = 1.054959511645151
$rowBF[0,3] = 1.0542130477284
$rowBF[0,4] = 1.064380378148236
$ws.Range("B3:F3").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.046568661151701
$rowIN[0,1] = 1.055276191699055
$rowIN[0,2] = 1.057516711613107
$rowIN[0,3] = 1.056772158639386
$rowIN[0,4] = 1.066913708975957
$rowIN[0,5] = 1.022193598356602
$ws.Range("I3:N3").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.051275972222913
$rowBF[0,2] = 1.055491982083952
$rowBF[0,3] = 1.054867324505244
$rowBF[0,4] = 1.065004122744296
$ws.Range("B4:F4").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.046762591504613
$rowIN[0,1] = 1.055734851351829
$rowIN[0,2] = 1.05792763277806
$rowIN[0,3] = 1.057304493731963
$rowIN[0,4] = 1.067416885983768
$rowIN[0,5] = 1.022349966238272
$ws.Range("I4:N4").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.051564650345293
$rowBF[0,2] = 1.055715888098839
$rowBF[0,3] = 1.055142526336449
$rowBF[0,4] = 1.065266477338492
$ws.Range("B5:F5").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.046843875253564
$rowIN[0,1] = 1.055927594751891
$rowIN[0,2] = 1.058100283887241
$rowIN[0,3] = 1.057528284890812
$rowIN[0,4] = 1.067628404706994
$rowIN[0,5] = 1.022415642653439
$ws.Range("I5:N5").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.051613124491273
$rowBF[0,2] = 1.055753486134812
$rowBF[0,3] = 1.055188742357484
$rowBF[0,4] = 1.065310535571592
$ws.Range("B6:F6").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.046857508803851
$rowIN[0,1] = 1.055959952684152
$rowIN[0,2] = 1.058129266882612
$rowIN[0,3] = 1.05756586025177
$rowIN[0,4] = 1.067663918625281
$rowIN[0,5] = 1.022426666449595
$ws.Range("I6:N6").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.051279829298172
$rowBF[0,2] = 1.055494973710104
$rowBF[0,3] = 1.054871001200862
$rowBF[0,4] = 1.06500762781996
$ws.Range("B7:F7").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.046763678583341
$rowIN[0,1] = 1.055737427102042
$rowIN[0,2] = 1.057929940144695
$rowIN[0,3] = 1.057307484050768
$rowIN[0,4] = 1.067419712374338
$rowIN[0,5] = 1.022350844048322
$ws.Range("I7:N7").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.049887262885853
$rowBF[0,2] = 1.054414906423549
$rowBF[0,3] = 1.053544117435405
$rowBF[0,4] = 1.063742648097504
$ws.Range("B8:F8").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.046369408626705
$rowIN[0,1] = 1.054806653022467
$rowIN[0,2] = 1.057095938199458
$rowIN[0,3] = 1.056227493521096
$rowIN[0,4] = 1.066398828244496
$rowIN[0,5] = 1.022033404776202
$ws.Range("I8:N8").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.047435470645317
$rowBF[0,2] = 1.052513516683916
$rowBF[0,3] = 1.05121061868048
$rowBF[0,4] = 1.061517878291053
$ws.Range("B9:F9").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.045666848930377
$rowIN[0,1] = 1.053164070505143
$rowIN[0,2] = 1.05562314981612
$rowIN[0,3] = 1.054324362976902
$rowIN[0,4] = 1.064599400383645
$rowIN[0,5] = 1.021472114063642
$ws.Range("I9:N9").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.045802345923186
$rowBF[0,2] = 1.051247191574641
$rowBF[0,3] = 1.049658099311652
$rowBF[0,4] = 1.060037618419168
$ws.Range("B10:F10").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.045193241249137
$rowIN[0,1] = 1.052067394771845
$rowIN[0,2] = 1.054639183858807
$rowIN[0,3] = 1.053055598499608
$rowIN[0,4] = 1.063399475974032
$rowIN[0,5] = 1.0210966412468
$ws.Range("I10:N10").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.045095516796122
$rowBF[0,2] = 1.050699167452419
$rowBF[0,3] = 1.04898659253446
$rowBF[0,4] = 1.059397349598437
$ws.Range("B11:F11").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.044986924868852
$rowIN[0,1] = 1.051592141199668
$rowIN[0,2] = 1.054212620648965
$rowIN[0,3] = 1.052506210989598
$rowIN[0,4] = 1.062879828694075
$rowIN[0,5] = 1.020933756627815
$ws.Range("I11:N11").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.044833017569734
$rowBF[0,2] = 1.050495652986029
$rowBF[0,3] = 1.048737277165227
$rowBF[0,4] = 1.059159629966478
$ws.Range("B12:F12").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.044910103466154
$rowIN[0,1] = 1.051415553116468
$rowIN[0,2] = 1.054054101350297
$rowIN[0,3] = 1.052302143635657
$rowIN[0,4] = 1.062686798119394
$rowIN[0,5] = 1.020873208778346
$ws.Range("I12:N12").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.044889322349428
$rowBF[0,2] = 1.050539305434934
$rowBF[0,3] = 1.048790751078279
$rowBF[0,4] = 1.059210616887349
$ws.Range("B13:F13").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.044926590346377
$rowIN[0,1] = 1.051453434469516
$rowIN[0,2] = 1.054088107654684
$rowIN[0,3] = 1.052345916788598
$rowIN[0,4] = 1.062728204298907
$rowIN[0,5] = 1.020886198542268
$ws.Range("I13:N13").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.045073817538196
$rowBF[0,2] = 1.050682343939548
$rowBF[0,3] = 1.048965981758669
$rowBF[0,4] = 1.059377697480254
$ws.Range("B14:F14").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.044980578587571
$rowIN[0,1] = 1.051577545550408
$rowIN[0,2] = 1.054199518916961
$rowIN[0,3] = 1.052489342724591
$rowIN[0,4] = 1.062863872937194
$rowIN[0,5] = 1.020928752645699
$ws.Range("I14:N14").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.045187497545162
$rowBF[0,2] = 1.050770480807999
$rowBF[0,3] = 1.049073962034222
$rowBF[0,4] = 1.059480655239673
$ws.Range("B15:F15").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.045013817841171
$rowIN[0,1] = 1.051654006837911
$rowIN[0,2] = 1.054268153179517
$rowIN[0,3] = 1.052577712135677
$rowIN[0,4] = 1.062947461491602
$rowIN[0,5] = 1.020954965642804
$ws.Range("I15:N15").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.045849262667985
$rowBF[0,2] = 1.051283568563399
$rowBF[0,3] = 1.049702680714485
$rowBF[0,4] = 1.060080125607954
$ws.Range("B16:F16").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.045206907651627
$rowIN[0,1] = 1.052098927676727
$rowIN[0,2] = 1.054667482957444
$rowIN[0,3] = 1.053092059473572
$rowIN[0,4] = 1.063433961767216
$rowIN[0,5] = 1.021107445006785
$ws.Range("I16:N16").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.046264457227631
$rowBF[0,2] = 1.051605496482523
$rowBF[0,3] = 1.050097258839147
$rowBF[0,4] = 1.060456343567817
$ws.Range("B17:F17").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.045327695609964
$rowIN[0,1] = 1.05237791163586
$rowIN[0,2] = 1.054917838475
$rowIN[0,3] = 1.05341469504107
$rowIN[0,4] = 1.063739111512986
$rowIN[0,5] = 1.021203010449346
$ws.Range("I17:N17").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.046506664526856
$rowBF[0,2] = 1.051793300788647
$rowBF[0,3] = 1.050327481381649
$rowBF[0,4] = 1.060675852013097
$ws.Range("B18:F18").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.045398029421156
$rowIN[0,1] = 1.052540601026555
$rowIN[0,2] = 1.055063818401586
$rowIN[0,3] = 1.053602882494896
$rowIN[0,4] = 1.063917093214622
$rowIN[0,5] = 1.021258722964459
$ws.Range("I18:N18").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.046589256277948
$rowBF[0,2] = 1.0518573421388
$rowBF[0,3] = 1.050405993487791
$rowBF[0,4] = 1.060750710047293
$ws.Range("B19:F19").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.045421991114975
$rowIN[0,1] = 1.052596067582929
$rowIN[0,2] = 1.055113585597306
$rowIN[0,3] = 1.05366704947571
$rowIN[0,4] = 1.063977779179151
$rowIN[0,5] = 1.021277714535515
$ws.Range("I19:N19").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.046219907540681
$rowBF[0,2] = 1.051570953639284
$rowBF[0,3] = 1.050054916911841
$rowBF[0,4] = 1.060415972020186
$ws.Range("B20:F20").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.0453147485915
$rowIN[0,1] = 1.052347983146901
$rowIN[0,2] = 1.054890982681727
$rowIN[0,3] = 1.053380079339044
$rowIN[0,4] = 1.063706372552483
$rowIN[0,5] = 1.021192760199758
$ws.Range("I20:N20").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.04501948696392
$rowBF[0,2] = 1.050640221388097
$rowBF[0,3] = 1.048914377592459
$rowBF[0,4] = 1.059328493515082
$ws.Range("B21:F21").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.044964685538068
$rowIN[0,1] = 1.051540999522861
$rowIN[0,2] = 1.054166713137499
$rowIN[0,3] = 1.052447107369487
$rowIN[0,4] = 1.062823922189631
$rowIN[0,5] = 1.0209162227705
$ws.Range("I21:N21").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.044265015897545
$rowBF[0,2] = 1.050055300046032
$rowBF[0,3] = 1.048197924258912
$rowBF[0,4] = 1.058645358735326
$ws.Range("B22:F22").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.044743509131497
$rowIN[0,1] = 1.051033282894629
$rowIN[0,2] = 1.053710903919306
$rowIN[0,3] = 1.051860508839436
$rowIN[0,4] = 1.062269030429761
$rowIN[0,5] = 1.020742091004517
$ws.Range("I22:N22").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.044664948712776
$rowBF[0,2] = 1.0503653524208
$rowBF[0,3] = 1.048577668102482
$rowBF[0,4] = 1.059007443777053
$ws.Range("B23:F23").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.044860861047414
$rowIN[0,1] = 1.051302464674379
$rowIN[0,2] = 1.053952577875165
$rowIN[0,3] = 1.052171475927431
$rowIN[0,4] = 1.062563194624186
$rowIN[0,5] = 1.020834426293681
$ws.Range("I23:N23").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.046240037529605
$rowBF[0,2] = 1.051586561973043
$rowBF[0,3] = 1.050074049184295
$rowBF[0,4] = 1.060434213978134
$ws.Range("B24:F24").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.045320599161857
$rowIN[0,1] = 1.052361506657783
$rowIN[0,2] = 1.054903117807638
$rowIN[0,3] = 1.053395720685988
$rowIN[0,4] = 1.06372116590104
$rowIN[0,5] = 1.02119739193652
$ws.Range("I24:N24").Value = $rowIN

$rowBF = New-Object 'object[,]' 1,5
$rowBF[0,0] = 1.02
$rowBF[0,1] = 1.048069070344481
$rowBF[0,2] = 1.053004850920849
$rowBF[0,3] = 1.051813331427423
$rowBF[0,4] = 1.062092522751373
$ws.Range("B25:F25").Value = $rowBF

$rowIN = New-Object 'object[,]' 1,6
$rowIN[0,0] = 1.045849401076917
$rowIN[0,1] = 1.053589004562598
$rowIN[0,2] = 1.056004274099353
$rowIN[0,3] = 1.054816371383592
$rowIN[0,4] = 1.065064651730353
$rowIN[0,5] = 1.021617447577732
$ws.Range("I25:N25").Value = $rowIN
